$wb = $excel.ActiveWorkbook

# Remember whatever sheet is active right now so we can restore focus to it
# at the end (selecting ranges on the other sheets below necessarily makes
# them the active sheet while the selection is made).
$originalActive = $wb.ActiveSheet

# The three fuel-list sheets all get the same two new rows appended:
#   row 12: "green hydrogen if"     -> 1
#   row 13: "low carbon hydrogen if" -> 1
$sheetPlan = @{
    "IFTQfS-PTC"   = "B33"
    "IFTQfS-ITC"   = "B14"
    "IFTQfS-loans" = "B13"
}

foreach ($name in $sheetPlan.Keys) {
    $ws = $wb.Worksheets.Item($name)

    # Duplicate the formatting (wrap-text style) of the last existing fuel
    # row onto the two new rows before filling in their values.
    $ws.Range("A11").Copy()
    $ws.Range("A12:A13").PasteSpecial(-4122)

    # The generic "hydrogen if" flag is superseded by the new, more
    # specific green/low-carbon hydrogen subscript rows, so it is turned
    # off here.
    $ws.Cells.Item(11, 2).Value = 0

    $ws.Cells.Item(12, 1).Value = "green hydrogen if"
    $ws.Cells.Item(12, 2).Value = 1

    $ws.Cells.Item(13, 1).Value = "low carbon hydrogen if"
    $ws.Cells.Item(13, 2).Value = 1

    # Leave the sheet's recorded selection where the author left it.
    $targetCell = $sheetPlan[$name]
    $ws.Range($targetCell).Select() | Out-Null
}

$originalActive.Activate() | Out-Null
